$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New file name strings (TC01 DNBSEQ-G400 -> TC04 Illumina HiSeq)
$neo4jFileName = 'TC04_CDS_Filter_InstrumentModel-Illumina HiSeq_Neo4jData.xlsx'
$webFileName   = 'TC04_CDS_Filter_InstrumentModel-Illumina HiSeq_WebData.xlsx'

# New Cypher queries (instrument_model filter DNBSEQ-G400 -> Illumina HiSeq)
$participantQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY `Participant ID`LIMIT 100
'@

$sampleQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@

$filesDetailQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@

$filesSummaryQuery = @'
MATCH (f:file)
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,f, s, collect(distinct samp.sample_id) as samp
RETURN
count(distinct s) AS Studies,
count(distinct p) AS Participants,
count(distinct samp) AS Samples,
count(distinct f) AS Files
'@

# Row 2 = ParticipantsTab
$ws.Range("B2").Value = $participantQuery
$ws.Range("C2").Value = $filesSummaryQuery
$ws.Range("D2").Value = $neo4jFileName
$ws.Range("E2").Value = $webFileName

# Row 3 = SamplesTab
$ws.Range("B3").Value = $sampleQuery
$ws.Range("C3").Value = $filesSummaryQuery
$ws.Range("D3").Value = $neo4jFileName
$ws.Range("E3").Value = $webFileName

# Row 4 = FilesTab
$ws.Range("B4").Value = $filesDetailQuery
$ws.Range("C4").Value = $filesSummaryQuery
$ws.Range("D4").Value = $neo4jFileName
$ws.Range("E4").Value = $webFileName

# Column D width change (78.85546875 -> 88.28515625)
$ws.Columns.Item(4).ColumnWidth = 87.5

# Selection change (B4 -> D4)
[void]$ws.Range("D4").Select()
